# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" for the 99a524a2-... row (row 7)
# on both locale sheets, reflecting a freshly generated handoff report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("D7").Value = "2016-03-11 06:03:55"
$wsDeDe.Range("D7").Value = "2016-03-11 06:04:04"
